$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 562, shifting existing rows 562..675 down to 563..676
$ws.Rows.Item(562).Insert()

# Populate the newly inserted row 562 with the new data record
$ws.Range("A562").Value = 3
$ws.Range("B562").Value = "Femacal de La Calera"
$ws.Range("C562").Value = "Coquimbo"
$ws.Range("D562").Value = 45258
$ws.Range("E562").Value = 5
$ws.Range("F562").Value = 100112027
$ws.Range("G562").Value = "Melón"
$ws.Range("H562").Value = "Tuna"
$ws.Range("I562").Value = "Primera"
$ws.Range("J562").Value = 50
$ws.Range("K562").Value = 17000
$ws.Range("L562").Value = 17000
$ws.Range("M562").Value = 17000
$ws.Range("N562").Value = "$/caja 12 unidades"
$ws.Range("O562").Value = "Provincia de Limarí"
$ws.Range("P562").Value = 1417
$ws.Range("Q562").Value = 12
$ws.Range("R562").Value = "Hortaliza"
